$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "GenValve: Std"
$ws.Range("L3").Value = "GenValve: Std"
$ws.Range("L4").Select()
